# Remove the two inline pictures (IMG_20251231_165044 and the
# "屏幕截图 2026-01-03 203917" screenshot) that were appended at the end of
# the document, while keeping their empty paragraphs (and the bookmark)
# intact.
$d = $word.ActiveDocument

for ($i = $d.InlineShapes.Count; $i -ge 1; $i--) {
    $shp = $d.InlineShapes.Item($i)
    $alt = $shp.AlternativeText
    if ($alt -eq "IMG_20251231_165044" -or $alt -eq "屏幕截图 2026-01-03 203917") {
        $shp.Delete()
    }
}
